# 20Questions.xlsx — add two new "20 Questions" decision trees (Adrian #4 /
# Giraffe in columns O/Q, and Adrian #5 / Penguin in columns S/U), plus
# finish off the existing Adrian #3 tree's leaf answer at J22 ("panda").
#
# Cells are written in the same order the strings were first typed so the
# resulting shared-string table ordering matches the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Finish Adrian #3 tree (panda guess)
$ws.Range('J22').Value = 'panda'

# --- Adrian #4: Giraffe tree (columns O / Q) ---
$ws.Range('O1').Value = 'Adrian #4'

$ws.Range('O2').Value = 'mammal'
$ws.Range('Q2').Value = 'yes'

$ws.Range('O3').Value = 'on land'
$ws.Range('Q3').Value = 'yes'

$ws.Range('O4').Value = '4 legs'
$ws.Range('Q4').Value = 'yes'

$ws.Range('O5').Value = 'ominvor'
$ws.Range('Q5').Value = 'no'

$ws.Range('O6').Value = 'carnivore'
$ws.Range('Q6').Value = 'yes'

$ws.Range('O7').Value = 'bigger than human'
$ws.Range('Q7').Value = 'YES'

$ws.Range('O8').Value = 'African'
$ws.Range('Q8').Value = 'yes'

$ws.Range('O9').Value = 'long neck'
$ws.Range('Q9').Value = 'yes'

$ws.Range('O10').Value = 'giraffe'
$ws.Range('Q10').Value = 'yes'

$ws.Range('O22').Value = 'Giraffe'

# --- Adrian #5: Penguin tree (columns S / U) ---
$ws.Range('S1').Value = 'Adrian #5'

$ws.Range('S3').Value = 'on land?'
$ws.Range('U3').Value = 'sometimes'

$ws.Range('S6').Value = '2 legs'

$ws.Range('S2').Value = 'reptile'
$ws.Range('U2').Value = 'yes'

$ws.Range('S4').Value = 'bigger than human'
$ws.Range('U4').Value = 'yes'

$ws.Range('S5').Value = '4 legs'
$ws.Range('U5').Value = 'no'

$ws.Range('U6').Value = 'no'

$ws.Range('S7').Value = 'carnivore'
$ws.Range('U7').Value = 'yes'

$ws.Range('S8').Value = 'crab?'
$ws.Range('U8').Value = 'no'

$ws.Range('S9').Value = 'vertibrate?'
$ws.Range('U9').Value = 'yes'

$ws.Range('S10').Value = 'live on beach?'
$ws.Range('U10').Value = 'no'

$ws.Range('S11').Value = 'near water?'
$ws.Range('U11').Value = 'yes'

$ws.Range('S12').Value = 'near river?'
$ws.Range('U12').Value = 'no'

$ws.Range('S13').Value = 'near lake?'
$ws.Range('U13').Value = 'no'

$ws.Range('S14').Value = 'cold temperatures?'
$ws.Range('U14').Value = 'yes'

$ws.Range('S15').Value = 'dangerous?'
$ws.Range('U15').Value = 'no'

$ws.Range('S16').Value = 'zoo?'
$ws.Range('U16').Value = 'yes'

$ws.Range('S17').Value = 'swim'
$ws.Range('U17').Value = 'yes'

$ws.Range('S22').Value = 'penguin'

# --- View state: scroll right a bit and leave the selection on W1 ---
$ws.Range('W1').Select()
$excel.ActiveWindow.Zoom = 100
